# Applies the "visibility map" update:
#  - Un-hide column C and column O, and rows 2 and 19-29 on the "Sheet1 (2)"
#    worksheet (they were previously hidden by the autofilter).
#  - Clear the stored colour autofilter condition (ShowAllData) while
#    keeping the AutoFilter itself in place.
#  - Add review comments ("Brian LeMaster: Complete") to C1, D1, E1 and F1
#    on that same worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1 (2)")
$ws.Activate()

# --- Un-hide previously filtered-out columns/rows -----------------------
$ws.Columns.Item(3).Hidden = $false    # column C
$ws.Columns.Item(15).Hidden = $false   # column O

$hiddenRows = @(2, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29)
foreach ($r in $hiddenRows) {
    $ws.Rows.Item($r).Hidden = $false
}

# --- Clear the autofilter's stored colour-filter condition --------------
$ws.ShowAllData()

# --- Add reviewer comments on the header cells ---------------------------
$commentText = "Brian LeMaster:`nComplete"
foreach ($col in @("C1", "D1", "E1", "F1")) {
    $rng = $ws.Range($col)
    $cmt = $rng.AddComment()
    $null = $cmt.Text($commentText)
    $cmt.Visible = $false
}
